# Updated cryptos list on Fri Jul 19 05:36:46 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "2" = @{ D = @{ V = "64.106.47"; T = "text" }; E = "  -0.81%  " }
    "3" = @{ D = @{ V = "3.419.83"; T = "text" }; E = "  +0.07%  " }
    "4" = @{ D = $null; E = "  +0.04%  " }
    "5" = @{ D = @{ V = "571.55"; T = "num" }; E = "  -0.07%  " }
    "6" = @{ D = @{ V = "160.93"; T = "num" }; E = "  +1.74%  " }
    "7" = @{ D = $null; E = "  +0.06%  " }
    "8" = @{ D = @{ V = "3.421.55"; T = "text" }; E = "  +0.03%  " }
    "9" = @{ D = @{ V = "0.552"; T = "num" }; E = "  -6.11%  " }
    "10" = @{ D = $null; E = "  +1.50%  " }
    "11" = @{ D = $null; E = "  -1.97%  " }
    "12" = @{ D = @{ V = "0.424"; T = "num" }; E = "  -2.90%  " }
    "13" = @{ D = @{ V = "4.012.14"; T = "text" }; E = "  +0.14%  " }
    "14" = @{ D = $null; E = "  +1.20%  " }
    "15" = @{ D = @{ V = "26.98"; T = "num" }; E = "  -2.43%  " }
    "16" = @{ D = $null; E = "  -6.87%  " }
    "17" = @{ D = @{ V = "64.148.01"; T = "text" }; E = "  -0.79%  " }
    "18" = @{ D = @{ V = "3.394.14"; T = "text" }; E = "  -1.04%  " }
    "19" = @{ D = $null; E = "  -4.26%  " }
    "20" = @{ D = $null; E = "  -1.97%  " }
    "21" = @{ D = @{ V = "376.15"; T = "num" }; E = "  -0.88%  " }
    "22" = @{ D = @{ V = "7.80"; T = "num" }; E = "  -2.41%  " }
    "23" = @{ D = @{ V = "0.998"; T = "num" }; E = "  -0.23%  " }
    "24" = @{ D = @{ V = "71.33"; T = "num" }; E = "  -0.66%  " }
    "25" = @{ D = $null; E = "  -5.31%  " }
    "26" = @{ D = @{ V = "0.0000115"; T = "num" }; E = "  -2.64%  " }
    "27" = @{ D = @{ V = "9.48"; T = "num" }; E = "  -4.42%  " }
    "28" = @{ D = $null; E = "  -0.01%  " }
    "29" = @{ D = $null; E = "  +0.08%  " }
    "30" = @{ D = @{ V = "6.00"; T = "num" }; E = "  -2.56%  " }
    "31" = @{ D = $null; E = "  -4.08%  " }
    "32" = @{ D = $null; E = "  +0.42%  " }
    "33" = @{ D = $null; E = "  -1.55%  " }
    "34" = @{ D = @{ V = "7.06"; T = "num" }; E = "  +0.96%  " }
    "35" = @{ D = $null; E = "  -3.77%  " }
    "36" = @{ D = @{ V = "159.57"; T = "num" }; E = "  -0.85%  " }
    "37" = @{ D = @{ V = "0.852"; T = "num" }; E = "  +10.56%  " }
    "38" = @{ D = @{ V = "1.81"; T = "num" }; E = "  -4.51%  " }
    "39" = @{ D = @{ V = "2.802.60"; T = "text" }; E = "  -2.75%  " }
    "40" = @{ D = $null; E = "  -3.38%  " }
    "41" = @{ D = @{ V = "42.90"; T = "num" }; E = "  -0.39%  " }
    "42" = @{ D = @{ V = "25.74"; T = "num" }; E = "  -1.84%  " }
    "43" = @{ D = @{ V = "6.43"; T = "num" }; E = "  -4.33%  " }
    "44" = @{ D = @{ V = "4.41"; T = "num" }; E = "  -2.43%  " }
    "45" = @{ D = @{ V = "25.91"; T = "num" }; E = "  +0.06%  " }
    "46" = @{ D = @{ V = "0.0305"; T = "num" }; E = "  -3.30%  " }
    "47" = @{ D = @{ V = "2.40"; T = "num" }; E = "  +7.54%  " }
    "48" = @{ D = @{ V = "334.71"; T = "num" }; E = "  +5.52%  " }
    "49" = @{ D = $null; E = "  -1.61%  " }
    "50" = @{ D = $null; E = "  -3.10%  " }
    "51" = @{ D = $null; E = "  -3.11%  " }
}

foreach ($rowNum in $updates.Keys) {
    $row = $updates[$rowNum]
    if ($null -ne $row.D) {
        $dcell = $ws.Range("D$rowNum")
        if ($row.D.T -eq "num") {
            # Force storage as text so values like "571.55" are not
            # reinterpreted as floating point numbers.
            $dcell.NumberFormat = "@"
        }
        $dcell.Value = $row.D.V
    }
    $ws.Range("E$rowNum").Value = $row.E
}
